$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 4.927909215762314
$ws.Range("D2").Value = 2.847003690747064
$ws.Range("E2").Value = 16.70553859340207
$ws.Range("F2").Value = 21.44205980522257
$ws.Range("G2").Value = 27.28524728152907
$ws.Range("H2").Value = 11.83963547182539
$ws.Range("N2").Value = 18.99769123883711
$ws.Range("O2").Value = 18.21424911005318
$ws.Range("C3").Value = 4.750711121090885
$ws.Range("D3").Value = 2.842352512309697
$ws.Range("E3").Value = 15.74098122386729
$ws.Range("F3").Value = 20.83462416852414
$ws.Range("G3").Value = 25.99870046402663
$ws.Range("H3").Value = 11.7393696284572
$ws.Range("N3").Value = 18.40031508502701
$ws.Range("O3").Value = 17.82091217572668
$ws.Range("C4").Value = 4.640096136786111
$ws.Range("D4").Value = 2.839518267372575
$ws.Range("E4").Value = 15.12306693296055
$ws.Range("F4").Value = 20.4598277646215
$ws.Range("G4").Value = 25.18411712988712
$ws.Range("H4").Value = 11.68066401629609
$ws.Range("N4").Value = 18.02485520896362
$ws.Range("O4").Value = 17.58096086381733
$ws.Range("C5").Value = 4.594643996562551
$ws.Range("D5").Value = 2.83836839667419
$ws.Range("E5").Value = 14.86507095140316
$ws.Range("F5").Value = 20.30690014922033
$ws.Range("G5").Value = 24.84649731098828
$ws.Range("H5").Value = 11.65748202894931
$ws.Range("N5").Value = 17.86990355188765
$ws.Range("O5").Value = 17.48372304092097
$ws.Range("C6").Value = 4.587076506565099
$ws.Range("D6").Value = 2.838177760592899
$ws.Range("E6").Value = 14.82186579043061
$ws.Range("F6").Value = 20.28150283021179
$ws.Range("G6").Value = 24.79011010989671
$ws.Range("H6").Value = 11.65367805690804
$ws.Range("N6").Value = 17.84406337566597
$ws.Range("O6").Value = 17.46761413378064
$ws.Range("C7").Value = 4.639484562621854
$ws.Range("D7").Value = 2.839502739525408
$ws.Range("E7").Value = 15.11961218647725
$ws.Range("F7").Value = 20.45776576031786
$ws.Range("G7").Value = 25.17958608262035
$ws.Range("H7").Value = 11.6803483472258
$ws.Range("N7").Value = 18.02277304767603
$ws.Range("O7").Value = 17.5796470742993
$ws.Range("C8").Value = 4.867234215169561
$ws.Range("D8").Value = 2.845395285005938
$ws.Range("E8").Value = 16.37843040157254
$ws.Range("F8").Value = 21.23315037590891
$ws.Range("G8").Value = 26.84705228545225
$ws.Range("H8").Value = 11.80448449313306
$ws.Range("N8").Value = 18.79364780656867
$ws.Range("O8").Value = 18.07839353513216
$ws.Range("C9").Value = 5.296259412424064
$ws.Range("D9").Value = 2.857137096684657
$ws.Range("E9").Value = 18.78811643531093
$ws.Range("F9").Value = 22.7282243724197
$ws.Range("G9").Value = 29.90193584892382
$ws.Range("H9").Value = 12.06954073893366
$ws.Range("N9").Value = 20.2273683202997
$ws.Range("O9").Value = 19.06226355741022
$ws.Range("C10").Value = 5.596955309575596
$ws.Range("D10").Value = 2.865899395333858
$ws.Range("E10").Value = 20.46965862704375
$ws.Range("F10").Value = 23.79757136625626
$ws.Range("G10").Value = 31.99356861721893
$ws.Range("H10").Value = 12.27596704518577
$ws.Range("N10").Value = 21.22223697909767
$ws.Range("O10").Value = 19.7802212038149
$ws.Range("C11").Value = 5.729934621946054
$ws.Range("D11").Value = 2.869919856340716
$ws.Range("E11").Value = 21.19208178177454
$ws.Range("F11").Value = 24.27519565339495
$ws.Range("G11").Value = 32.90844354024738
$ws.Range("H11").Value = 12.37207576034201
$ws.Range("N11").Value = 21.66018057919901
$ws.Range("O11").Value = 20.10410721387743
$ws.Range("C12").Value = 5.779695868351379
$ws.Range("D12").Value = 2.871447636775248
$ws.Range("E12").Value = 21.45956983581143
$ws.Range("F12").Value = 24.45460272838332
$ws.Range("G12").Value = 33.24938747018248
$ws.Range("H12").Value = 12.40875721936211
$ws.Range("N12").Value = 21.82377585682186
$ws.Range("O12").Value = 20.22623451673249
$ws.Range("C13").Value = 5.769006056086128
$ws.Range("D13").Value = 2.871118359485898
$ws.Range("E13").Value = 21.40223123933002
$ws.Range("F13").Value = 24.41603192273011
$ws.Range("G13").Value = 33.17620698650827
$ws.Range("H13").Value = 12.40084491739543
$ws.Range("N13").Value = 21.78864458690801
$ws.Range("O13").Value = 20.19995739655366
$ws.Range("C14").Value = 5.734040688071927
$ws.Range("D14").Value = 2.87004543903785
$ws.Range("E14").Value = 21.21420978851934
$ws.Range("F14").Value = 24.28998584363616
$ws.Range("G14").Value = 32.93660461711558
$ws.Range("H14").Value = 12.3750879827921
$ws.Range("N14").Value = 21.67368539489659
$ws.Range("O14").Value = 20.11416586801264
$ws.Range("C15").Value = 5.712544578162985
$ws.Range("D15").Value = 2.869388947936345
$ws.Range("E15").Value = 21.09825076922159
$ws.Range("F15").Value = 24.21258366221619
$ws.Range("G15").Value = 32.78911880636149
$ws.Range("H15").Value = 12.35934762366352
$ws.Range("N15").Value = 21.60297336126124
$ws.Range("O15").Value = 20.06154445023211
$ws.Range("C16").Value = 5.58818396741828
$ws.Range("D16").Value = 2.865637373057486
$ws.Range("E16").Value = 20.42159408957463
$ws.Range("F16").Value = 23.76616327202328
$ws.Range("G16").Value = 31.933021644086
$ws.Range("H16").Value = 12.26972793351934
$ws.Range("N16").Value = 21.19330956972086
$ws.Range("O16").Value = 19.75898847577349
$ws.Range("C17").Value = 5.510881612523735
$ws.Range("D17").Value = 2.863344984699761
$ws.Range("E17").Value = 19.99562288820602
$ws.Range("F17").Value = 23.48989760279738
$ws.Range("G17").Value = 31.39827571243351
$ws.Range("H17").Value = 12.21529267921034
$ws.Range("N17").Value = 20.93814219015166
$ws.Range("O17").Value = 19.57258892858086
$ws.Range("C18").Value = 5.466063290532993
$ws.Range("D18").Value = 2.862029683714533
$ws.Range("E18").Value = 19.74661129124327
$ws.Range("F18").Value = 23.33017627434371
$ws.Range("G18").Value = 31.08726844452466
$ws.Range("H18").Value = 12.18419226319307
$ws.Range("N18").Value = 20.79000725568362
$ws.Range("O18").Value = 19.46512847085719
$ws.Range("C19").Value = 5.450828966686399
$ws.Range("D19").Value = 2.861584887426505
$ws.Range("E19").Value = 19.66161140627431
$ws.Range("F19").Value = 23.27596249551343
$ws.Range("G19").Value = 30.98138443994493
$ws.Range("H19").Value = 12.17369904556637
$ws.Range("N19").Value = 20.73962067985785
$ws.Range("O19").Value = 19.42870566065305
$ws.Range("C20").Value = 5.519147792434978
$ws.Range("D20").Value = 2.863588677445738
$ws.Range("E20").Value = 20.04138240899327
$ws.Range("F20").Value = 23.51939291460245
$ws.Range("G20").Value = 31.45555771327277
$ws.Range("H20").Value = 12.22106595754394
$ws.Range("N20").Value = 20.96544799483449
$ws.Range("O20").Value = 19.59245815574117
$ws.Range("C21").Value = 5.744327365249977
$ws.Range("D21").Value = 2.870360434484887
$ws.Range("E21").Value = 21.26960091134836
$ws.Range("F21").Value = 24.32704965068447
$ws.Range("G21").Value = 33.00713256012475
$ws.Range("H21").Value = 12.3826458645136
$ws.Range("N21").Value = 21.70751365554064
$ws.Range("O21").Value = 20.13938006887667
$ws.Range("C22").Value = 5.888008340429265
$ws.Range("D22").Value = 2.87481734652353
$ws.Range("E22").Value = 22.03690264483966
$ws.Range("F22").Value = 24.84632078162225
$ws.Range("G22").Value = 33.98904006068562
$ws.Range("H22").Value = 12.48990748496821
$ws.Range("N22").Value = 22.17935961385674
$ws.Range("O22").Value = 20.49373176747093
$ws.Range("C23").Value = 5.811656343170425
$ws.Range("D23").Value = 2.872435633780768
$ws.Range("E23").Value = 21.63060736680948
$ws.Range("F23").Value = 24.57001917168679
$ws.Range("G23").Value = 33.46798535939462
$ws.Range("H23").Value = 12.43251800643284
$ws.Range("N23").Value = 21.92877110911181
$ws.Range("O23").Value = 20.30493171273619
$ws.Range("C24").Value = 5.515411820516509
$ws.Range("D24").Value = 2.863478495811163
$ws.Range("E24").Value = 20.02070734994363
$ws.Range("F24").Value = 23.50606084954937
$ws.Range("G24").Value = 31.42967163699083
$ws.Range("H24").Value = 12.21845524834599
$ws.Range("N24").Value = 20.95310750188673
$ws.Range("O24").Value = 19.58347619283895
$ws.Range("C25").Value = 5.182497424840949
$ws.Range("D25").Value = 2.853938411188361
$ws.Range("E25").Value = 18.13092776440016
$ws.Range("F25").Value = 22.32795688141718
$ws.Range("G25").Value = 29.10099330607924
$ws.Range("H25").Value = 11.9956708575648
$ws.Range("N25").Value = 19.84905939529497
$ws.Range("O25").Value = 18.79637008422408
